# Update the cryptos list with fresh prices / volume(1h) percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)
# Price values are stored as text (they use '.' as a thousands separator,
# and some carry significant trailing zeros), so force the cell to Text
# format before writing, otherwise Excel would coerce them into numbers.

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "62.556.48"
$ws.Cells.Item(2, 5).Value = "  +0.81%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.024.51"
$ws.Cells.Item(3, 5).Value = "  +1.25%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.00"

# Row 5 - BNB
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "596.98"
$ws.Cells.Item(5, 5).Value = "  +1.55%  "

# Row 6 - Solana
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "150.84"
$ws.Cells.Item(6, 5).Value = "  +5.45%  "

# Row 7 - USDC
$ws.Cells.Item(7, 5).Value = "  -0.11%  "

# Row 8 - LidoStakedEther
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "3.023.75"
$ws.Cells.Item(8, 5).Value = "  +1.34%  "

# Row 9 - XRP
$ws.Cells.Item(9, 5).Value = "  -0.56%  "

# Row 10 - Toncoin
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "6.47"
$ws.Cells.Item(10, 5).Value = "  +12.61%  "

# Row 11 - Dogecoin
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.151"
$ws.Cells.Item(11, 5).Value = "  +3.65%  "

# Row 12 - Cardano
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.461"
$ws.Cells.Item(12, 5).Value = "  -0.25%  "

# Row 13 - ShibaInu
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000234"
$ws.Cells.Item(13, 5).Value = "  +3.10%  "

# Row 14 - Avalanche
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "34.65"
$ws.Cells.Item(14, 5).Value = "  +1.44%  "

# Row 15 - TRON
$ws.Cells.Item(15, 5).Value = "  +2.34%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "3.523.36"
$ws.Cells.Item(16, 5).Value = "  +0.86%  "

# Row 17 - Polkadot
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "7.04"
$ws.Cells.Item(17, 5).Value = "  +0.45%  "

# Row 18 - WrappedBTC
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "62.502.65"
$ws.Cells.Item(18, 5).Value = "  +0.64%  "

# Row 19 - WrappedEther
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "3.020.28"
$ws.Cells.Item(19, 5).Value = "  +0.81%  "

# Row 20 - BitcoinCash
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "449.90"
$ws.Cells.Item(20, 5).Value = "  -0.66%  "

# Row 21 - Chainlink
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "14.26"
$ws.Cells.Item(21, 5).Value = "  +2.85%  "

# Row 22 - Polygon
$ws.Cells.Item(22, 5).Value = "  +1.57%  "

# Row 23 - Uniswap
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "7.45"
$ws.Cells.Item(23, 5).Value = "  +1.43%  "

# Row 24 - Litecoin
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "82.45"
$ws.Cells.Item(24, 5).Value = "  +1.14%  "

# Row 25 - Fetch.AI
$ws.Cells.Item(25, 5).Value = "  +4.02%  "

# Row 26 - RenderToken
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "10.82"
$ws.Cells.Item(26, 5).Value = "  +12.80%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "12.05"
$ws.Cells.Item(27, 5).Value = "  -0.73%  "

# Row 28 - Dai
$ws.Cells.Item(28, 5).Value = "  +0.16%  "

# Row 29 - PancakeSwap
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.72"
$ws.Cells.Item(29, 5).Value = "  +3.15%  "

# Row 30 - FirstDigitalUSD
$ws.Cells.Item(30, 5).Value = "  -0.21%  "

# Row 31 - NEARProtocol
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.23"
$ws.Cells.Item(31, 5).Value = "  +5.12%  "

# Row 32 - ImmutableX
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "2.15"
$ws.Cells.Item(32, 5).Value = "  +3.79%  "

# Row 33 - EthereumClassic
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "27.65"
$ws.Cells.Item(33, 5).Value = "  +0.65%  "

# Row 34 - Hedera
$ws.Cells.Item(34, 5).Value = "  +2.71%  "

# Row 35 - PEPE
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.0₃0857"
$ws.Cells.Item(35, 5).Value = "  +8.44%  "

# Row 36 - Mantle
$ws.Cells.Item(36, 5).Value = "  +1.07%  "

# Row 37 - Filecoin
$ws.Cells.Item(37, 5).Value = "  +2.98%  "

# Row 38 - dogwifhat
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "3.06"
$ws.Cells.Item(38, 5).Value = "  +8.19%  "

# Row 39 - Stacks
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.09"
$ws.Cells.Item(39, 5).Value = "  +0.14%  "

# Row 40 - OKB
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "50.16"
$ws.Cells.Item(40, 5).Value = "  +0.18%  "

# Row 41 - Cosmos
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "9.06"
$ws.Cells.Item(41, 5).Value = "  -0.40%  "

# Row 42 - Kaspa
$ws.Cells.Item(42, 5).Value = "  +2.37%  "

# Row 43 & 44 - Arweave and TheGraph swap order (TheGraph now ranks above Arweave)
$ws.Cells.Item(43, 2).Value = "TheGraph"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.286"
$ws.Cells.Item(43, 5).Value = "  +8.15%  "

$ws.Cells.Item(44, 2).Value = "Arweave"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "40.95"
$ws.Cells.Item(44, 5).Value = "  +11.33%  "

# Row 45 - Bittensor
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "394.18"
$ws.Cells.Item(45, 5).Value = "  +0.87%  "

# Row 46 - VeChain
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0354"
$ws.Cells.Item(46, 5).Value = "  +0.02%  "

# Row 47 - Maker
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.755.58"
$ws.Cells.Item(47, 5).Value = "  +0.97%  "

# Row 48 - Monero
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "134.83"
$ws.Cells.Item(48, 5).Value = "  +4.22%  "

# Row 49 - USDe
$ws.Cells.Item(49, 5).Value = "  +0.07%  "

# Row 50 - ThetaToken
$ws.Cells.Item(50, 5).Value = "  +0.91%  "

# Row 51 - Stellar
$ws.Cells.Item(51, 5).Value = "  -0.34%  "
